$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update row 40 values
$ws.Range("B40").Value = 314
$ws.Range("C40").Value = 394
$ws.Range("D40").Value = 390
$ws.Range("E40").Value = 66

# Update row 41 values
$ws.Range("B41").Value = 426
$ws.Range("C41").Value = 458
$ws.Range("D41").Value = 166
$ws.Range("E41").Value = 66

# Update the selected cell in the sheet view from E43 to C43
$ws.Range("C43").Select()
